$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in remarks for the ICs (F27:F30) -- ATTINY first so new shared
# strings land in the same append order as the target workbook.
$ws.Range("F30").Value = "ATTINY 84A-SSU"
$ws.Range("F27").Value = "SMD HC 595"
$ws.Range("F28").Value = "SMD 40106"
$ws.Range("F29").Value = "SMD 4028"

# Update the document version string (B5)
$ws.Range("B5").Value = "Document Version 03/04/2023"

# Update the active selection to D6
$ws.Range("D6").Select()
